# Refresh crypto Price (D) and 1h-change (E) columns with latest scraped values.
# D-column values are numeric-looking text (e.g. '0.999', '59.864.62') that must
# stay stored as text, exactly like the rest of the sheet's inline/shared strings --
# so each write temporarily forces Text format, then restores the default 'Normal'
# style afterwards (Excel keeps the already-entered value as text even after the
# number format is reset, matching the workbook's original no-explicit-style cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '59.864.62'; E = '  +1.12%  ' }
    @{ Row = 3; D = '2.300.49'; E = '  -0.67%  ' }
    @{ Row = 4; D = '0.999'; E = '  -0.06%  ' }
    @{ Row = 5; D = '541.23'; E = '  +0.12%  ' }
    @{ Row = 6; D = '129.29'; E = '  -2.50%  ' }
    @{ Row = 7; D = '0.999'; E = '  -0.11%  ' }
    @{ Row = 8; D = '0.574'; E = '  -2.27%  ' }
    @{ Row = 9; D = '2.300.12'; E = '  -0.49%  ' }
    @{ Row = 10; D = '0.101'; E = '  -0.48%  ' }
    @{ Row = 11; D = '5.51'; E = '  +0.48%  ' }
    @{ Row = 12; D = $null; E = '  -0.23%  ' }
    @{ Row = 13; D = '0.331'; E = '  -0.97%  ' }
    @{ Row = 14; D = '23.24'; E = '  -3.04%  ' }
    @{ Row = 15; D = '59.747.85'; E = '  +1.19%  ' }
    @{ Row = 16; D = '2.706.48'; E = '  -0.72%  ' }
    @{ Row = 17; D = '0.0000132'; E = '  -1.30%  ' }
    @{ Row = 18; D = '2.282.64'; E = '  -0.98%  ' }
    @{ Row = 19; D = '10.44'; E = '  -1.94%  ' }
    @{ Row = 20; D = '4.07'; E = '  -2.75%  ' }
    @{ Row = 21; D = '311.39'; E = '  -0.89%  ' }
    @{ Row = 22; D = '6.55'; E = '  -0.77%  ' }
    @{ Row = 23; D = $null; E = '  -0.21%  ' }
    @{ Row = 24; D = '5.68'; E = '  -0.48%  ' }
    @{ Row = 25; D = '63.65'; E = '  +1.51%  ' }
    @{ Row = 26; D = '0.169'; E = '  -2.24%  ' }
    @{ Row = 27; D = '1.00'; E = '  +0.09%  ' }
    @{ Row = 28; D = '7.72'; E = '  -3.19%  ' }
    @{ Row = 29; D = '1.34'; E = '  +2.33%  ' }
    @{ Row = 30; D = '170.51'; E = '  +0.02%  ' }
    @{ Row = 31; D = '1.17'; E = '  -0.46%  ' }
    @{ Row = 32; D = '1.71'; E = '  -1.10%  ' }
    @{ Row = 33; D = '0.0₃0723'; E = '  -2.47%  ' }
    @{ Row = 34; D = '5.81'; E = '  -1.52%  ' }
    @{ Row = 35; D = '1.36'; E = '  +1.72%  ' }
    @{ Row = 36; D = '0.377'; E = '  -1.94%  ' }
    @{ Row = 37; D = $null; E = '  +0.01%  ' }
    @{ Row = 38; D = '17.63'; E = '  -1.36%  ' }
    @{ Row = 39; D = $null; E = '  -0.11%  ' }
    @{ Row = 40; D = '3.99'; E = '  -2.87%  ' }
    @{ Row = 41; D = '317.35'; E = '  +4.29%  ' }
    @{ Row = 42; D = '37.82'; E = '  -1.77%  ' }
    @{ Row = 43; D = $null; E = '  -1.40%  ' }
    @{ Row = 44; D = '136.04'; E = '  -3.67%  ' }
    @{ Row = 45; D = '3.41'; E = '  -1.41%  ' }
    @{ Row = 46; D = '0.0935'; E = '  -2.48%  ' }
    @{ Row = 47; D = '0.562'; E = '  +0.85%  ' }
    @{ Row = 48; D = '18.69'; E = '  +0.92%  ' }
    @{ Row = 49; D = '0.0488'; E = '  -1.57%  ' }
    @{ Row = 50; D = '0.0₆0222'; E = '  +19.47%  ' }
    @{ Row = 51; D = '0.0211'; E = '  -0.49%  ' }
)

foreach ($u in $updates) {
    if ($u.D -ne $null) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
